# issue #5: stock data output to json file
#
# 1. Clean up a handful of shared-string typos in the "股票" (stock) sheet:
#    - drop a stray inner space in several "XXX股份有限公司" names
#    - drop the full-width comma in a few quantity/total numbers that
#      were stored as text (keep them as text, just fix the digits)
# 2. Add a new "property_category" column (all rows = "stock") between
#    the existing "total" and "date" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- 1a. Fix company-name cells that had a stray inner space -----------
$nameFixes = @{
    7  = "台灣土地開發股份有限公司"
    8  = "景岳生物科技股份有限公司"
    13 = "利奇機械工業股份有限公司"
    16 = "台灣玻璃工業股份有限公司"
    20 = "中華開發金融控股股份有限公司"
    21 = "玉山金融控股股份有限公司"
    23 = "寶島光學科技股份有限公司"
    25 = "上福全球科技股份有限公H]"
    26 = "富旺國際開發股份有限公司"
}
foreach ($row in $nameFixes.Keys) {
    $ws.Cells.Item($row, 2).Value = $nameFixes[$row]
}

# --- 1b. Fix numeric-looking text cells that used a full-width comma ---
# (column D row 6, column G rows 11/18/28) - these must stay text cells,
# so force a text number format before writing the digits.
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "35000"

$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "65000"

$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "20000"

$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "10000"

# --- 2. Insert the new "property_category" column (H), pushing the old
#        date / legislator_name / legislator_id columns one to the right.
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
